$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 07:22"

# Row 61 - Kazajistan
$ws.Range("B61").Value = 3205
$ws.Range("C61").Value = 67
$ws.Range("E61").Value = 2361

# Row 62 - Tailandia
$ws.Range("B62").Value = 2954
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 2684
$ws.Range("E62").Value = 216

# Row 64 - Hungria
$ws.Range("B64").Value = 2775
$ws.Range("C64").Value = 48
$ws.Range("D64").Value = 581
$ws.Range("E64").Value = 1882
$ws.Range("F64").Value = 54
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 312

# Row 83 - Eslovenia
$ws.Range("D83").Value = 1091
$ws.Range("E83").Value = 238
$ws.Range("F83").Value = 22
